$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 3246
$ws.Range("I63").Value = 3246
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3246
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2622
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 2992.739
$ws.Range("I64").Value = 2782
$ws.Range("J64").Value = 3051.2778
$ws.Range("K64").Value = 2782
$ws.Range("L64").Value = 3051.2778
$ws.Range("M64").Value = -2534
$ws.Range("N64").Value = -3547.2778

$ws.Range("H66").Value = 3246
$ws.Range("I66").Value = 3246
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9738
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6618
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 2992.739
$ws.Range("I67").Value = 2782
$ws.Range("J67").Value = 3051.2778
$ws.Range("K67").Value = 2782
$ws.Range("L67").Value = 3051.2778
$ws.Range("M67").Value = -1924
$ws.Range("N67").Value = -4767.2778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 115.22222
$ws.Range("I4").Value = 106.42857
$ws.Range("J4").Value = 146
$ws.Range("K4").Value = 106.42857
$ws.Range("L4").Value = 146
$ws.Range("M4").Value = 9.571430000000007
$ws.Range("N4").Value = -378

$ws.Range("H63").Value = 3177.4092
$ws.Range("I63").Value = 2783.8333
$ws.Range("J63").Value = 3325
$ws.Range("K63").Value = 2783.8333
$ws.Range("L63").Value = 3325
$ws.Range("M63").Value = -2097.8333
$ws.Range("N63").Value = -4697

$ws.Range("H66").Value = 3177.4092
$ws.Range("I66").Value = 2783.8333
$ws.Range("J66").Value = 3325
$ws.Range("K66").Value = 13919.1665
$ws.Range("L66").Value = 16625
$ws.Range("M66").Value = -10487.1665
$ws.Range("N66").Value = -23489

$ws.Range("H96").Value = 23315.6
$ws.Range("J96").Value = 23315.6
$ws.Range("L96").Value = 23315.6
$ws.Range("N96").Value = -28807.6

$ws.Range("H113").Value = 30080
$ws.Range("J113").Value = 30080
$ws.Range("L113").Value = 30080
$ws.Range("N113").Value = -38758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 29207.875
$ws.Range("J106").Value = 29207.875
$ws.Range("L106").Value = 29207.875
$ws.Range("N106").Value = -31731.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17858578
$ws.Range("I31").Value = 1331.25
$ws.Range("K31").Value = 1331.25
$ws.Range("M31").Value = -1036.25

$ws.Range("H34").Value = 17858578
$ws.Range("I34").Value = 1331.25
$ws.Range("K34").Value = 1331.25
$ws.Range("M34").Value = -1129.25

$ws.Range("H62").Value = 2424.7778
$ws.Range("I62").Value = 2266.6667
$ws.Range("J62").Value = 2582.889
$ws.Range("K62").Value = 2266.6667
$ws.Range("L62").Value = 2582.889
$ws.Range("M62").Value = -1642.6667
$ws.Range("N62").Value = -3830.889

$ws.Range("H65").Value = 2424.7778
$ws.Range("I65").Value = 2266.6667
$ws.Range("J65").Value = 2582.889
$ws.Range("K65").Value = 11333.3335
$ws.Range("L65").Value = 12914.445
$ws.Range("M65").Value = -8213.333500000001
$ws.Range("N65").Value = -19154.445

$ws.Range("H68").Value = 23333.334
$ws.Range("J68").Value = 23333.334
$ws.Range("L68").Value = 23333.334
$ws.Range("N68").Value = -24831.334

$ws.Range("H71").Value = 23333.334
$ws.Range("J71").Value = 23333.334
$ws.Range("L71").Value = 70000.00199999999
$ws.Range("N71").Value = -77488.00199999999

$ws.Range("H107").Value = 1132.9642
$ws.Range("I107").Value = 444.94446
$ws.Range("J107").Value = 2371.4
$ws.Range("K107").Value = 444.94446
$ws.Range("L107").Value = 2371.4
$ws.Range("M107").Value = 1475.05554
$ws.Range("N107").Value = -6211.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 3654.2
$ws.Range("I82").Value = 2241
$ws.Range("J82").Value = 4596.3335
$ws.Range("K82").Value = 6723
$ws.Range("L82").Value = 13789.0005
$ws.Range("M82").Value = -6317
$ws.Range("N82").Value = -14601.0005

$ws.Range("H85").Value = 3654.2
$ws.Range("I85").Value = 2241
$ws.Range("J85").Value = 4596.3335
$ws.Range("K85").Value = 6723
$ws.Range("L85").Value = 13789.0005
$ws.Range("M85").Value = -5319
$ws.Range("N85").Value = -16597.0005

$ws.Range("H104").Value = 2306.1538
$ws.Range("I104").Value = 1000
$ws.Range("J104").Value = 2415
$ws.Range("K104").Value = 3000
$ws.Range("L104").Value = 7245
$ws.Range("M104").Value = -379
$ws.Range("N104").Value = -12487

$ws.Range("H131").Value = 928.04
$ws.Range("I131").Value = 832
$ws.Range("J131").Value = 933.0947
$ws.Range("K131").Value = 2496
$ws.Range("L131").Value = 2799.2841
$ws.Range("M131").Value = 2544
$ws.Range("N131").Value = -12879.2841

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 24000
$ws.Range("J103").Value = 24000
$ws.Range("L103").Value = 24000
$ws.Range("N103").Value = -26344

$ws.Range("H107").Value = 692.7632
$ws.Range("I107").Value = 627.125
$ws.Range("J107").Value = 805.2857
$ws.Range("K107").Value = 627.125
$ws.Range("L107").Value = 805.2857
$ws.Range("M107").Value = 1292.875
$ws.Range("N107").Value = -4645.2857

$ws.Range("H112").Value = 28858.4
$ws.Range("J112").Value = 28858.4
$ws.Range("L112").Value = 28858.4
$ws.Range("N112").Value = -31074.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 22860
$ws.Range("I68").Value = 51500
$ws.Range("J68").Value = 3766.6667
$ws.Range("K68").Value = 51500
$ws.Range("L68").Value = 3766.6667
$ws.Range("M68").Value = -50751
$ws.Range("N68").Value = -5264.6667

$ws.Range("H71").Value = 22860
$ws.Range("I71").Value = 51500
$ws.Range("J71").Value = 3766.6667
$ws.Range("K71").Value = 257500
$ws.Range("L71").Value = 18833.3335
$ws.Range("M71").Value = -253756
$ws.Range("N71").Value = -26321.3335

$ws.Range("H104").Value = 21894.285
$ws.Range("J104").Value = 21894.285
$ws.Range("L104").Value = 21894.285
$ws.Range("N104").Value = -28882.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8706.076999999999
$ws.Range("I107").Value = 13585.125
$ws.Range("K107").Value = 40755.375
$ws.Range("M107").Value = -38835.375
